$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12; this shifts existing rows 12..87 down to 13..88
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44749
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100114007
$ws.Range("G12").Value = "Jengibre"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 170
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 12412
$ws.Range("N12").Value = "$/caja 13 kilos"
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 955
$ws.Range("Q12").Value = 13
$ws.Range("R12").Value = "Hortaliza"

# Ensure the date cell keeps the same date number format as the rest of column D
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat
